$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "76.257.26"
$ws.Cells.Item(2, 5).Value = "  -0.47%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.077.35"
$ws.Cells.Item(3, 5).Value = "  +5.05%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "198.74"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.38%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "617.96"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.90%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.34%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.210"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +5.90%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "3.076.25"
$ws.Cells.Item(10, 5).Value = "  +4.99%  "

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.446"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.29%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.04%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +7.57%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.641.73"
$ws.Cells.Item(14, 5).Value = "  +4.87%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +3.47%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "ShibaInu"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000195"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.75%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "76.164.39"
$ws.Cells.Item(17, 5).Value = "  -0.48%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.072.04"
$ws.Cells.Item(18, 5).Value = "  +4.63%  "

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "13.55"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.63%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +2.99%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +14.61%  "

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "382.29"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.30%  "

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +4.54%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +1.10%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "NEARProtocol"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +7.83%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "WrappedeETH"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(26, 4).Value = "3.234.55"
$ws.Cells.Item(26, 5).Value = "  +4.78%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "72.41"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.55%  "

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.04%  "

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "10.07"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.96%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.98%  "

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.31%  "

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "8.31"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.24%  "

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +4.60%  "

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "502.21"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.10%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +5.93%  "

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.03%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +13.38%  "

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "20.82"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.53%  "

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "162.22"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.06%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "194.73"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +8.50%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.50%  "

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.28%  "

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -7.30%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.02%  "

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +24.03%  "

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "5.17"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +5.37%  "

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.26"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +6.64%  "

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.66"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.82%  "

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +6.32%  "

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "40.62"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.46%  "

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.598"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.77%  "
